# Actualización 11 de Mayo - Mañana
$wb = $excel.ActiveWorkbook

# --- Hoja "1er Parcial" ---
$ws1 = $wb.Worksheets.Item("1er Parcial")

# Fila 8 (4ARHM)
$ws1.Cells.Item(8, 5).Value = 18
$ws1.Cells.Item(8, 6).Value = 6
$ws1.Cells.Item(8, 7).Value = 75
$ws1.Cells.Item(8, 8).Value = 25
$ws1.Cells.Item(8, 9).Value = 7.5
$ws1.Cells.Item(8, 10).Value = 2
$ws1.Cells.Item(8, 11).Value = 8.33

# Fila 13 (4BLCM)
$ws1.Cells.Item(13, 5).Value = 28
$ws1.Cells.Item(13, 6).Value = 8
$ws1.Cells.Item(13, 7).Value = 77.78
$ws1.Cells.Item(13, 8).Value = 22.22
$ws1.Cells.Item(13, 10).Value = 8
$ws1.Cells.Item(13, 11).Value = 22.22

# --- Hoja "2o Parcial" ---
$ws2 = $wb.Worksheets.Item("2o Parcial")

# Fila 6 (4AEM)
$ws2.Cells.Item(6, 5).Value = 28
$ws2.Cells.Item(6, 6).Value = 11
$ws2.Cells.Item(6, 7).Value = 71.79
$ws2.Cells.Item(6, 8).Value = 28.21
$ws2.Cells.Item(6, 9).Value = 8.2
$ws2.Cells.Item(6, 10).Value = 10
$ws2.Cells.Item(6, 11).Value = 25.64

# Fila 7 (4ALCM)
$ws2.Cells.Item(7, 5).Value = 29
$ws2.Cells.Item(7, 6).Value = 8
$ws2.Cells.Item(7, 7).Value = 78.38
$ws2.Cells.Item(7, 8).Value = 21.62
$ws2.Cells.Item(7, 9).Value = 8.8
$ws2.Cells.Item(7, 10).Value = 8
$ws2.Cells.Item(7, 11).Value = 21.62

# Fila 8 (4ARHM)
$ws2.Cells.Item(8, 5).Value = 18
$ws2.Cells.Item(8, 6).Value = 6
$ws2.Cells.Item(8, 7).Value = 75
$ws2.Cells.Item(8, 8).Value = 25
$ws2.Cells.Item(8, 9).Value = 8.7
$ws2.Cells.Item(8, 10).Value = 5
$ws2.Cells.Item(8, 11).Value = 20.83

# Fila 9 (4BEM)
$ws2.Cells.Item(9, 5).Value = 17
$ws2.Cells.Item(9, 6).Value = 20
$ws2.Cells.Item(9, 7).Value = 45.95
$ws2.Cells.Item(9, 8).Value = 54.05
$ws2.Cells.Item(9, 9).Value = 7.7
$ws2.Cells.Item(9, 10).Value = 16
$ws2.Cells.Item(9, 11).Value = 43.24

# Fila 10 (4ALCV)
$ws2.Cells.Item(10, 5).Value = 21
$ws2.Cells.Item(10, 6).Value = 12
$ws2.Cells.Item(10, 7).Value = 63.64
$ws2.Cells.Item(10, 8).Value = 36.36
$ws2.Cells.Item(10, 9).Value = 7.9
$ws2.Cells.Item(10, 10).Value = 11
$ws2.Cells.Item(10, 11).Value = 33.33

# Fila 11 (4ASV)
$ws2.Cells.Item(11, 5).Value = 18
$ws2.Cells.Item(11, 6).Value = 16
$ws2.Cells.Item(11, 7).Value = 52.94
$ws2.Cells.Item(11, 8).Value = 47.06
$ws2.Cells.Item(11, 9).Value = 8.2
$ws2.Cells.Item(11, 10).Value = 15
$ws2.Cells.Item(11, 11).Value = 44.12

# Fila 12 (4APM)
$ws2.Cells.Item(12, 5).Value = 12
$ws2.Cells.Item(12, 6).Value = 16
$ws2.Cells.Item(12, 7).Value = 42.86
$ws2.Cells.Item(12, 8).Value = 57.14
$ws2.Cells.Item(12, 9).Value = 7.1
$ws2.Cells.Item(12, 10).Value = 16
$ws2.Cells.Item(12, 11).Value = 57.14

# Fila 13 (4BLCM)
$ws2.Cells.Item(13, 5).Value = 25
$ws2.Cells.Item(13, 6).Value = 11
$ws2.Cells.Item(13, 7).Value = 69.44
$ws2.Cells.Item(13, 8).Value = 30.56
$ws2.Cells.Item(13, 9).Value = 8.4
$ws2.Cells.Item(13, 10).Value = 11
$ws2.Cells.Item(13, 11).Value = 30.56

# --- Hoja "3er Parcial" ---
$ws3 = $wb.Worksheets.Item("3er Parcial")

# Fila 6 (4AEM)
$ws3.Cells.Item(6, 5).Value = 32
$ws3.Cells.Item(6, 6).Value = 7
$ws3.Cells.Item(6, 7).Value = 82.05
$ws3.Cells.Item(6, 8).Value = 17.95
$ws3.Cells.Item(6, 9).Value = 7.4

# Fila 7 (4ALCM)
$ws3.Cells.Item(7, 9).Value = 8.1

# Fila 8 (4ARHM)
$ws3.Cells.Item(8, 5).Value = 20
$ws3.Cells.Item(8, 6).Value = 4
$ws3.Cells.Item(8, 7).Value = 83.33
$ws3.Cells.Item(8, 8).Value = 16.67
$ws3.Cells.Item(8, 9).Value = 8.1
$ws3.Cells.Item(8, 10).Value = 2
$ws3.Cells.Item(8, 11).Value = 8.33

# Fila 9 (4BEM)
$ws3.Cells.Item(9, 9).Value = 6.7

# Fila 10 (4ALCV)
$ws3.Cells.Item(10, 5).Value = 23
$ws3.Cells.Item(10, 6).Value = 10
$ws3.Cells.Item(10, 7).Value = 69.7
$ws3.Cells.Item(10, 8).Value = 30.3
$ws3.Cells.Item(10, 9).Value = 7.4

# Fila 11 (4ASV)
$ws3.Cells.Item(11, 5).Value = 22
$ws3.Cells.Item(11, 6).Value = 12
$ws3.Cells.Item(11, 7).Value = 64.71
$ws3.Cells.Item(11, 8).Value = 35.29
$ws3.Cells.Item(11, 9).Value = 7.9

# Fila 13 (4BLCM)
$ws3.Cells.Item(13, 5).Value = 28
$ws3.Cells.Item(13, 6).Value = 8
$ws3.Cells.Item(13, 7).Value = 77.78
$ws3.Cells.Item(13, 8).Value = 22.22
$ws3.Cells.Item(13, 9).Value = 8.2
$ws3.Cells.Item(13, 10).Value = 8
$ws3.Cells.Item(13, 11).Value = 22.22
